$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Kmeans example completed, log updated
#
# Row 25: new "K mean python notebook example" task, marked Done / Medium,
#          1 hour, notes "completed", next steps "research project ideas".
# Row 26: new "Find new dataset and perform basic k mean model" task,
#          also Done / Medium, 1 hour.
# ---------------------------------------------------------------------------

# Seed row 25 / row 26 date cells from row 24 (same date, same date style)
# so they pick up the existing date number-format style instead of minting
# a brand new style entry.
$ws.Range("A24").Copy($ws.Range("A25"))
$ws.Range("A24").Copy($ws.Range("A26"))

# Row 25
$ws.Range("A25").Value = 45707
$ws.Range("B25").Value = "K mean python notebook example"
$ws.Range("C25").Value = "Done"
$ws.Range("D25").Value = "Medium"
$ws.Range("E25").Value = 1
$ws.Range("F25").Value = "completed"
$ws.Range("G25").Value = "research project ideas"

# Row 26
$ws.Range("A26").Value = 45707
$ws.Range("B26").Value = "Find new dataset and perform basic k mean model"
$ws.Range("C26").Value = "Done"
$ws.Range("D26").Value = "Medium"
$ws.Range("E26").Value = 1

# Update the window view state to reflect where the user ended up after
# logging the entries (scrolled down a bit further, selection on F28).
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("F28").Select()
